$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:G2").Copy()
$ws.Range("A23:G23").PasteSpecial(-4122)

$ws.Range("A23").Value = "Windows10"
$ws.Range("B23").Value = "PC"
$ws.Range("C23").Value = "Ryzen 9 7950x"
$ws.Range("D23").Value = "4.5"
$ws.Range("E23").Value = "DDR4-3200 32GB"
$ws.Range("F23").Value = 61
$ws.Range("G23").Value = 45737

$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1

$ws.Range("A24").Select()
